$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = " 56"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = " 44"
$ws.Range("B3").Value = 94
$ws.Range("C3").Value = 1.37
$ws.Range("D3").Value = 309.8
$ws.Range("E3").Value = 1.68
$ws.Range("F3").Value = 197.4
$ws.Range("G3").Value = 0.73
$ws.Range("H3").Value = 1.16
$ws.Range("I3").Value = 0.05
$ws.Range("J3").Value = 0.2
$ws.Range("K3").Value = 0.1
$ws.Range("L3").Value = 109
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 5
$ws.Range("O3").Value = 19
$ws.Range("P3").Value = 9
